$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.403.42'
$ws.Range('E2').Value = '  +2.34%  '
$ws.Range('D3').Value = '1.574.85'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +1.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  +1.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.95'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.76'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0592'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '1.799.91'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').Value = '1.571.58'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '28.400.42'
$ws.Range('E17').Value = '  +2.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('D21').Value = '0.0₃0693'
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('E23').Value = '  -4.18%  '
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '150.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('E28').Value = '  -1.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.104'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('E31').Value = '  -1.86%  '
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.19'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = '1.396.35'
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.55'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.75%  '
$ws.Range('E37').Value = '  -3.82%  '
$ws.Range('E38').Value = '  +2.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.57'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.06%  '
$ws.Range('E40').Value = '  -1.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.533'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.793'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.980'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '62.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.12%  '
$ws.Range('D48').Value = '1.711.49'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.60'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0101'
$ws.Range('E51').Value = '  +0.80%  '
